$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# "Antigua BD" fix: periodo mora 2508 -> 2509 for Novedad de Ingreso (MARILUZ, row 20),
# EDWIN's matching period row (row 21) and Novedad de Retiro (STEEL POLO, row 22) all
# shared the same "2508" label - update them all to "2509".
$ws.Range("E20").Value = "2509"
$ws.Range("E21").Value = "2509"
$ws.Range("E22").Value = "2509"

# Center the "Periodo Mora" values for the whole data table (rows 16-22).
$ws.Range("E16:E22").HorizontalAlignment = -4108
